$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the Price column as Text first so numeric-looking strings
# (e.g. "542.99", "0.999") are stored as literal text, matching the
# original inlineStr cells instead of being parsed into numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "60.988.82"
$ws.Range("E2").Value = "  +3.99%  "
$ws.Range("D3").Value = "2.550.07"
$ws.Range("E3").Value = "  +3.29%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "542.99"
$ws.Range("E5").Value = "  +1.97%  "
$ws.Range("D6").Value = "146.72"
$ws.Range("E6").Value = "  +1.85%  "
$ws.Range("D7").Value = "0.995"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("D9").Value = "2.590.88"
$ws.Range("E9").Value = "  +4.00%  "
$ws.Range("E10").Value = "  +2.83%  "
$ws.Range("E11").Value = "  +1.50%  "
$ws.Range("D12").Value = "5.54"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("E13").Value = "  +4.26%  "
$ws.Range("D14").Value = "3.004.97"
$ws.Range("E14").Value = "  +3.62%  "
$ws.Range("D15").Value = "24.53"
$ws.Range("E15").Value = "  +3.32%  "
$ws.Range("D16").Value = "60.787.01"
$ws.Range("E16").Value = "  +3.84%  "
$ws.Range("D17").Value = "0.0000144"
$ws.Range("E17").Value = "  +5.46%  "
$ws.Range("D18").Value = "2.565.47"
$ws.Range("E18").Value = "  +3.63%  "
$ws.Range("D19").Value = "11.40"
$ws.Range("E19").Value = "  +1.27%  "
$ws.Range("E20").Value = "  +2.10%  "
$ws.Range("D21").Value = "329.28"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("D22").Value = "5.99"
$ws.Range("E22").Value = "  +4.72%  "
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").Value = "63.42"
$ws.Range("E24").Value = "  +4.56%  "
$ws.Range("D25").Value = "0.444"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("E26").Value = "  +4.89%  "
$ws.Range("D27").Value = "0.993"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("D28").Value = "8.10"
$ws.Range("E28").Value = "  +5.50%  "
$ws.Range("D29").Value = "7.21"
$ws.Range("E29").Value = "  +4.10%  "
$ws.Range("D30").Value = "0.0₃0815"
$ws.Range("E30").Value = "  +5.70%  "
$ws.Range("E31").Value = "  +2.66%  "
$ws.Range("D32").Value = "1.22"
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("D33").Value = "164.47"
$ws.Range("E33").Value = "  +4.16%  "
$ws.Range("D34").Value = "1.48"
$ws.Range("E34").Value = "  +6.13%  "
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").Value = "18.91"
$ws.Range("E36").Value = "  +2.41%  "
$ws.Range("D37").Value = "4.51"
$ws.Range("E37").Value = "  +3.04%  "
$ws.Range("D38").Value = "1.66"
$ws.Range("E38").Value = "  +4.27%  "
$ws.Range("D39").Value = "5.74"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").Value = "309.01"
$ws.Range("E40").Value = "  +1.47%  "
$ws.Range("D41").Value = "37.16"
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("D42").Value = "0.849"
$ws.Range("E42").Value = "  +5.23%  "
$ws.Range("E43").Value = "  +2.36%  "
$ws.Range("E44").Value = "  +3.99%  "
$ws.Range("D45").Value = "0.997"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "10.88"
$ws.Range("E46").Value = "  +1.03%  "
$ws.Range("D47").Value = "126.90"
$ws.Range("E47").Value = "  +2.22%  "
$ws.Range("E48").Value = "  +3.90%  "
$ws.Range("D49").Value = "0.0942"
$ws.Range("E49").Value = "  +2.17%  "
$ws.Range("E50").Value = "  +1.66%  "
$ws.Range("D51").Value = "0.0232"
$ws.Range("E51").Value = "  +2.18%  "

# Restore the default (unformatted) style on the Price column so the
# cells keep the workbook's original appearance/style index.
$priceRange.Style = "Normal"
